# BL Audit Form update: "added 2nd-10-24 report"
# - shift the report date forward a day (01.10.2024 -> 02.10.2024)
# - shift the payment note forward a day (02.10.2024 payment -> 03.10.2024 payment)
# - update the day's stock/balance figures (dependent formulas recalc automatically)
# - move the active selection to F13 with the view scrolled back to the top

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B1: "01.10.2024" -> "02.10.2024" -------------------------------------
# A plain dd.mm.yyyy-looking string typed straight into Range.Value gets
# auto-coerced to a date serial (and picks up a new date number-format),
# which would change the cell's type/style. Build the text via a formula
# in a scratch cell instead, then paste-special just the value back onto
# B1 so the literal string lands with the original "General" style intact.
$scratch = $ws.Range("H1")
$scratch.Formula = "=""02.10.2024"""
$scratch.Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$scratch.ClearContents()

# --- F32: "02.10.2024 payment " -> "03.10.2024 payment " -------------------
# Trailing text keeps this from being read as a date, so a direct value
# assignment is safe here.
$ws.Range("F32").Value = "03.10.2024 payment "

# --- updated figures for the day (formulas in column E recalc on their own)
$ws.Range("C9").Value = 338933
$ws.Range("C10").Value = 560
$ws.Range("C11").Value = 17490
$ws.Range("C14").Value = 15
$ws.Range("E20").Value = 33794
$ws.Range("E21").Value = 10249
$ws.Range("E25").Value = 24750
$ws.Range("E32").Value = 118000

# --- restore the view to the top and move the selection to F13 -------------
$ws.Range("F13").Select() | Out-Null
